# Auto-generated files on 2025-09-23
# Update the HotStock_Top20 rankings (columns A:C, rows 2-21) with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "张江高科"
$ws.Range("B2").Value = "上海建工"
$ws.Range("C2").Value = "天普股份"

$ws.Range("A3").Value = "和而泰"
$ws.Range("B3").Value = "立讯精密"
$ws.Range("C3").Value = "凯美特气"

$ws.Range("A4").Value = "凯美特气"
$ws.Range("B4").Value = "山子高科"
$ws.Range("C4").Value = "张江高科"

$ws.Range("A5").Value = "立讯精密"
$ws.Range("B5").Value = "张江高科"
$ws.Range("C5").Value = "立讯精密"

$ws.Range("A6").Value = "天普股份"
$ws.Range("B6").Value = "沃尔核材"
$ws.Range("C6").Value = "山子高科"

$ws.Range("A7").Value = "长川科技"
$ws.Range("B7").Value = "和而泰"
$ws.Range("C7").Value = "沃尔核材"

$ws.Range("A8").Value = "山子高科"
$ws.Range("B8").Value = "凯美特气"
$ws.Range("C8").Value = "和而泰"

$ws.Range("A9").Value = "上海建工"
$ws.Range("B9").Value = "首开股份"
$ws.Range("C9").Value = "杭电股份"

$ws.Range("A10").Value = "工业富联"
$ws.Range("B10").Value = "大众公用"
$ws.Range("C10").Value = "上海建工"

$ws.Range("A11").Value = "沃尔核材"
$ws.Range("B11").Value = "工业富联"
$ws.Range("C11").Value = "三花智控"

$ws.Range("A12").Value = "立昂微"
$ws.Range("B12").Value = "福龙马"
$ws.Range("C12").Value = "长川科技"

$ws.Range("A13").Value = "三花智控"
$ws.Range("B13").Value = "长川科技"
$ws.Range("C13").Value = "万向钱潮"

$ws.Range("A14").Value = "杭电股份"
$ws.Range("B14").Value = "东华软件"
$ws.Range("C14").Value = "福龙马"

$ws.Range("A15").Value = "首开股份"
$ws.Range("B15").Value = "立昂微"
$ws.Range("C15").Value = "首开股份"

$ws.Range("A16").Value = "大众公用"
$ws.Range("B16").Value = "大洋电机"
$ws.Range("C16").Value = "大洋电机"

$ws.Range("A17").Value = "大洋电机"
$ws.Range("B17").Value = "万向钱潮"
$ws.Range("C17").Value = "卧龙电驱"

$ws.Range("A18").Value = "福龙马"
$ws.Range("B18").Value = "三花智控"
$ws.Range("C18").Value = "东华软件"

$ws.Range("A19").Value = "波长光电"
$ws.Range("B19").Value = "永鼎股份"
$ws.Range("C19").Value = "大众公用"

$ws.Range("A20").Value = "万向钱潮"
$ws.Range("B20").Value = "天普股份"
$ws.Range("C20").Value = "海立股份"

$ws.Range("A21").Value = "东华软件"
$ws.Range("B21").Value = "先导智能"
$ws.Range("C21").Value = "先导智能"
